$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix the blank "LIS" 8/22 result (E50) to "N/A" (reuses the existing
#        shared string / keeps the existing s="1" style already on that cell) ---
$ws.Range("E50").Value = "N/A"

# --- 2. Re-style the Sept-6 batch (rows 57:65, column E) so it uses the same
#        plain boxed-border style as the rest of the sheet instead of the
#        one-off "fill + left/right border" style ---
$src1 = $ws.Range("E2")
$dst1 = $ws.Range("E57:E65")
$src1.Copy()
$dst1.PasteSpecial(-4122)

# --- 3. Append the new Sept-18/19 rows (66:76), formatted like the existing
#        data rows, then fill in their values ---
$fmtSrc = $ws.Range("A51:E51")
$newRows = $ws.Range("A66:E76")
$fmtSrc.Copy()
$newRows.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A66").Value = "RCS"
$ws.Range("B66").Value = "ES0923B0578"
$ws.Range("C66").Value = 45188
$ws.Range("D66").Value = 0.4465277777777778
$ws.Range("E66").Value = 3

$ws.Range("A67").Value = "RD22"
$ws.Range("B67").Value = "ES0923B0580"
$ws.Range("C67").Value = 45188
$ws.Range("D67").Value = 0.31805555555555554
$ws.Range("E67").Value = 3

$ws.Range("A68").Value = "I80"
$ws.Range("B68").Value = "ES0923B0582"
$ws.Range("C68").Value = 45188
$ws.Range("D68").Value = 0.30833333333333335
$ws.Range("E68").Value = 2

$ws.Range("A69").Value = "LIS"
$ws.Range("B69").Value = "ES0923B0583"
$ws.Range("C69").Value = 45188
$ws.Range("D69").Value = 0.38541666666666669
$ws.Range("E69").Value = 2

$ws.Range("A70").Value = "STTD"
$ws.Range("B70").Value = "ES0923B0584"
$ws.Range("C70").Value = 45188
$ws.Range("D70").Value = 0.47430555555555554
$ws.Range("E70").Value = 4

$ws.Range("A71").Value = "BL5"
$ws.Range("B71").Value = "ES0923B0585"
$ws.Range("C71").Value = 45188
$ws.Range("D71").Value = 0.39930555555555558
$ws.Range("E71").Value = 4

$ws.Range("A72").Value = "PRS"
$ws.Range("B72").Value = "ES0923B0586"
$ws.Range("C72").Value = 45188
$ws.Range("D72").Value = 0.375
$ws.Range("E72").Value = 2

$ws.Range("A73").Value = "LIB"
$ws.Range("B73").Value = "ES0923B0587"
$ws.Range("C73").Value = 45188
$ws.Range("D73").Value = 0.35138888888888892
$ws.Range("E73").Value = 1

$ws.Range("A74").Value = "RYI"
$ws.Range("B74").Value = "ES0923B0588"
$ws.Range("C74").Value = 45188
$ws.Range("D74").Value = 0.32222222222222224
$ws.Range("E74").Value = 2

$ws.Range("A75").Value = "RVB"
$ws.Range("B75").Value = "ES0923B0589"
$ws.Range("C75").Value = 45188
$ws.Range("D75").Value = 0.28958333333333336
$ws.Range("E75").Value = 3

$ws.Range("A76").Value = "SHR"
$ws.Range("B76").Value = "ES0923B0590"
$ws.Range("C76").Value = 45187
$ws.Range("D76").Value = 0.36388888888888887
$ws.Range("E76").Value = 1

# --- 4. Keep the selection on the last data row, same as the source file ---
$ws.Range("E65").Select() | Out-Null
